{"js": "// The author renamed the helper function they define in the write-up from\n// `generate_constraints()` to `generate_con()` (the later sentence that\n// *calls* the function already said `generate_con()`, so only the\n// definition mention needs to change).\nconst body = context.document.body;\nconst results = body.search(\"generate_constraints()\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"generate_con()\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# The author renamed the helper function they define in the write-up from\n# `generate_constraints()` to `generate_con()` (the later sentence that\n# *calls* the function already said `generate_con()`, so only the\n# definition mention needs to change).\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"generate_constraints()\"\n$find.Replacement.Text = \"generate_con()\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.Forward = $true\n$find.Wrap = 1 # wdFindContinue\n\n# wdReplace: 0 = wdReplaceNone, 1 = wdReplaceOne, 2 = wdReplaceAll\n$find.Execute($null, $null, $null, $null, $null, $null, $true, $null, $null, $null, 2)\n"}
